# Edit script for DataWorkbookRubella.xlsx
# Commit message: "closer population, better births and fixed 0-14 deaths"
#
# Changes applied:
#   - On the "Initial_conditions" sheet, column D (rows 2-109) holds a
#     constant value that is used (among C and E) to compute the "Check"
#     column M (SUM(C:E)-G). The value is changed from 10 to 3 for every
#     data row (2 through 109).
#   - The active cell selection on that sheet moves from G8 to T85.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial_conditions")
$ws.Activate()

# Update D2:D109 from 10 to 3 for all data rows.
$ws.Range("D2:D109").Value = 3

# Move the selection to match the saved view state.
$ws.Range("T85").Select()
